$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.913.96"
$ws.Range("E2").Value = "  +1.19%  "

$ws.Range("D3").Value = "1.638.38"
$ws.Range("E3").Value = "  +0.55%  "

$ws.Range("E4").Value = "  +0.13%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.53"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("E6").Value = "  +0.53%  "

$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.61"
$ws.Range("E8").Value = "  +0.98%  "

$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0876"
$ws.Range("E11").Value = "  -0.16%  "

$ws.Range("E12").Value = "  +0.63%  "

$ws.Range("D13").Value = "1.643.31"
$ws.Range("E13").Value = "  +0.87%  "

$ws.Range("E14").Value = "  +4.02%  "

$ws.Range("E15").Value = "  +1.07%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.99"

$ws.Range("D17").Value = "27.907.78"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "231.89"
$ws.Range("E18").Value = "  +0.43%  "

$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").Value = "  +0.87%  "

$ws.Range("E20").Value = "  +0.56%  "

$ws.Range("E21").Value = "  +0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.76"
$ws.Range("E22").Value = "  +1.99%  "

$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("E24").Value = "  -3.73%  "

$ws.Range("E25").Value = "  +1.57%  "

$ws.Range("E26").Value = "  +0.29%  "

$ws.Range("E27").Value = "  +1.30%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("E29").Value = "  +0.22%  "

$ws.Range("E30").Value = "  +0.86%  "

$ws.Range("E31").Value = "  -0.07%  "

$ws.Range("E32").Value = "  +1.69%  "

$ws.Range("E33").Value = "  +1.50%  "

$ws.Range("D34").Value = "1.409.45"
$ws.Range("E34").Value = "  -4.22%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.58"
$ws.Range("E35").Value = "  +1.25%  "

$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.890"
$ws.Range("E37").Value = "  +1.09%  "

$ws.Range("E38").Value = "  +0.11%  "

$ws.Range("E39").Value = "  -0.09%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.916"
$ws.Range("E40").Value = "  -2.73%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.03"
$ws.Range("E41").Value = "  -0.62%  "

$ws.Range("E42").Value = "  +0.07%  "

$ws.Range("E43").Value = "  +5.19%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "66.17"
$ws.Range("E44").Value = "  -2.47%  "

$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").Value = "1.779.73"
$ws.Range("E47").Value = "  +0.76%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "88.12"
$ws.Range("E48").Value = "  +0.61%  "

$ws.Range("E49").Value = "  +0.85%  "

$ws.Range("E50").Value = "  +0.40%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.64"
$ws.Range("E51").Value = "  -1.10%  "
